# Weekly price update: insert the latest week's record for
# "Vega Modelo de Temuco - Achicoria" above the existing row 149,
# pushing the previously-recorded rows (149-168) down to (150-169).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 149 (shifts 149:168 -> 150:169).
$ws.Rows.Item(149).Insert()

# Populate the newly inserted row with the new weekly observation.
$ws.Range("A149").Value = 10
$ws.Range("B149").Value = "Vega Modelo de Temuco"
$ws.Range("C149").Value = "La Araucanía"
$ws.Range("D149").Value = 45275
$ws.Range("E149").Value = 9
$ws.Range("F149").Value = 100112010
$ws.Range("G149").Value = "Achicoria"
$ws.Range("H149").Value = "Sin especificar"
$ws.Range("I149").Value = "Primera"
$ws.Range("J149").Value = 35
$ws.Range("K149").Value = 10000
$ws.Range("L149").Value = 10000
$ws.Range("M149").Value = 10000
$ws.Range("N149").Value = "$/caja 18 unidades"
$ws.Range("O149").Value = "Región Metropolitana"
$ws.Range("P149").Value = 556
$ws.Range("Q149").Value = 18
$ws.Range("R149").Value = "Hortaliza"
